$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Row 2 (Beta) - columns F through N
$ws.Range("F2").Value = 8.651971006319151
$ws.Range("G2").Value = 7.405160994755866
$ws.Range("H2").Value = 9.838323441170894
$ws.Range("I2").Value = 1.600710017157438
$ws.Range("J2").Value = 1.26706293279254
$ws.Range("K2").Value = 1.877295042496091
$ws.Range("L2").Value = 0.1259838013027723
$ws.Range("M2").Value = 0.1043275539874963
$ws.Range("N2").Value = 0.1443046060300968

# Update Row 3 (Gamma) - columns F through N
$ws.Range("F3").Value = 0.1059789297132179
$ws.Range("G3").Value = 0.00007240790001257524
$ws.Range("H3").Value = 0.3022823134360977
$ws.Range("I3").Value = 0.0908312168431834
$ws.Range("J3").Value = 0.00006205386501834856
$ws.Range("K3").Value = 0.2607515211087289
$ws.Range("L3").Value = 0.10701031826124
$ws.Range("M3").Value = 0.0000726335646425335
$ws.Range("N3").Value = 0.3054671810784932

# Add new Row 4 (Beta + Gamma) - copy style of A3 into A4 first
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 12.00687180793019
$ws.Range("D4").Value = 1.974477778970852
$ws.Range("E4").Value = 0.1537386519519979
$ws.Range("F4").Value = 8.75794993603237
$ws.Range("G4").Value = 7.405233402655877
$ws.Range("H4").Value = 10.14060575460699
$ws.Range("I4").Value = 1.691541234000621
$ws.Range("J4").Value = 1.267124986657559
$ws.Range("K4").Value = 2.13804656360482
$ws.Range("L4").Value = 0.2329941195640122
$ws.Range("M4").Value = 0.1044001875521388
$ws.Range("N4").Value = 0.4497717871085901

Write-Output "done"
